# Update the example model (2_species_1_reaction) so that:
#  a) the model changes mass during the simulation (species type molecular
#     weights differ, and the rate law "Equation" becomes the literal
#     string "100"), and
#  b) selection/view state reflects the sheets touched while making the
#     edits (Species types, then Rate laws as the final active sheet).

$wb = $excel.ActiveWorkbook

# --- Species types: change molecular weights (column F) so the model's
#     mass changes during the simulation -------------------------------
$species = $wb.Worksheets.Item("Species types")
$species.Range("F2").Value = 1000
$species.Range("F3").Value = 1
$species.Range("F5").Select()

# --- Rate laws: the rate law's "Equation" (column C, row 2) becomes the
#     quoted literal string "100" and is formatted with a (scientific)
#     number format; the column is also widened to fit ------------------
$rateLaws = $wb.Worksheets.Item("Rate laws")
$rateLaws.Range("C2").Value = '"100"'
$rateLaws.Range("C2").NumberFormat = "0.00E+00"
$rateLaws.Columns.Item(3).ColumnWidth = 12.666666666666666

# Rate laws ends up as the active sheet/tab, with C3 selected
$rateLaws.Range("C3").Select()
